$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: fill in E2/F2 (observations) and G2 (status) ---
$ws.Range("E2").Value = "Todas as câmeras estão funcionando."
$ws.Range("F2").Value = "Troca de ballun e p4."
$ws.Range("G2").Value = "Concluido"

# --- Row 5: fill in E5/F5 (observations) and G5 (status); F5 becomes wrap-text ---
$ws.Range("E5").Value = "Comunicação foi passada para o DDNS, falta eu subir aqui (estava de almoço)."
$ws.Range("F5").Value = '"necessario mandar o vendedor ao local para orcamento, procurar por Thales o mais rapido possivel."'
$ws.Range("F5").WrapText = $true
$ws.Range("G5").Value = "Concluido"
$ws.Rows.Item(5).RowHeight = 30

# --- Row 8: fill in E8 and G8 ---
$ws.Range("E8").Value = "Marcos esteve no local, acredito que tenha sido resolvido."
$ws.Range("G8").Value = "Em andamento"

# --- Row 9: fill in E9 and G9 ---
$ws.Range("E9").Value = "Marcos esteve no local, acredito que tenha sido resolvido."
$ws.Range("G9").Value = "Em andamento"

# --- Row 12: fill in E12/F12 and G12 ---
$ws.Range("E12").Value = "Necessário trocar a central do cliente."
$ws.Range("F12").Value = "Giovani disse que é necessário trocar a centarl do cliente, Active 20."
$ws.Range("G12").Value = "Concluido"

# --- Row 13: fill in E13 and flip G13 from "Em andamento" to "Concluido" ---
$ws.Range("E13").Value = "Aparentemente central estava travada. Técnico disse que pediu o Marcos substituição do chip GPRS."
$ws.Range("G13").Value = "Concluido"

# --- Insert a brand-new row at 14 (Estivas / totem incident), pushing
#     everything below down by one. ---
$ws.Rows.Item(14).Insert()

# Fill the new row's values first (while default column styles are active),
# then paste the formatting (borders / quote-prefix numfmt) from the row that
# is now two rows below (old "Fábio" row, now row 16) so the new row matches
# the look of the rest of the table exactly.
$ws.Range("B14").Value = "0003"
$ws.Range("C14").Value = "Estivas"
$ws.Range("D14").Value = "Batida de carro no totem."
$ws.Range("E14").Value = "Pedro foi prestar a devida avalição no totem e reparo, foi repassado ao Marcos."
$ws.Range("G14").Value = "Concluido"

$ws.Range("A16:H16").Copy()
$ws.Range("A14:H14").PasteSpecial(-4122)

# --- Update selection to match the target workbook ---
$ws.Range("F2").Select() | Out-Null
